$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.420.75'
$ws.Range('E2').Value = '  -4.07%  '
$ws.Range('D3').Value = '3.565.03'
$ws.Range('E3').Value = '  -4.54%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '587.19'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -4.30%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '184.85'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('D7').Value = '3.559.35'
$ws.Range('E7').Value = '  -4.65%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.614'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -4.10%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -7.11%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.146'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -9.82%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '53.09'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -7.42%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000261'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -11.38%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.82'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -8.12%  '
$ws.Range('D15').Value = '4.134.32'
$ws.Range('E15').Value = '  -4.43%  '
$ws.Range('D16').Value = '3.559.51'
$ws.Range('E16').Value = '  -4.75%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.126'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('E18').Value = '  -5.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.25'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -6.58%  '
$ws.Range('D20').Value = '66.311.14'
$ws.Range('E20').Value = '  -3.91%  '
$ws.Range('E21').Value = '  -7.33%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '396.10'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -4.45%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.36'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -5.76%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '85.99'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.83%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.37'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +3.58%  '
$ws.Range('E26').Value = '  -5.34%  '
$ws.Range('E27').Value = '  -3.50%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.05'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.54'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -6.83%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.96'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -7.75%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '31.13'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -6.61%  '
$ws.Range('E32').Value = '  -3.44%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '12.18'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.81%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '619.99'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.54%  '
$ws.Range('E35').Value = '  -7.56%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '63.08'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.30%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '41.35'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -7.29%  '
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.401'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('D40').Value = '0.0₃0767'
$ws.Range('E40').Value = '  -9.98%  '
$ws.Range('E41').Value = '  -6.81%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = '3.009.19'
$ws.Range('E43').Value = '  +6.24%  '
$ws.Range('E44').Value = '  -8.28%  '
$ws.Range('E45').Value = '  -4.58%  '
$ws.Range('E46').Value = '  -7.85%  '
$ws.Range('E47').Value = '  -7.80%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.11'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.56'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -7.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '137.67'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.38%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.73'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.80%  '
